$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 8: Sunday entry ---
$ws.Range("A8").Value = "Sunday"
$ws.Range("B8").Value = 45774
$ws.Range("C8").Value = 0.54166666666666663
$ws.Range("D8").Value = 0.70833333333333337
$ws.Range("E8").Formula = "=D8-C8"
$ws.Range("F8").Value = "Suffering on camera movement and character rotation, working on zoom"

# Copy number formats from the row above so styles dedupe against the
# existing date/time cell styles instead of creating new duplicate ones.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("C7:E7").Copy()
$ws.Range("C8:E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Totals column ---
$ws.Range("H1").Value = "Total time spent"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Interior.Color = 65535

$ws.Range("I1").Formula = "=SUM(E:E)"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").NumberFormat = "[$-F400]h:mm:ss AM/PM"

# --- Column widths ---
$ws.Columns.Item(6).ColumnWidth = 66.7109375
$ws.Columns.Item(8).ColumnWidth = 15.5703125
$ws.Columns.Item(9).ColumnWidth = 13.5703125

# --- Selection ---
$ws.Range("J9").Select()
